$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaps")

# Clear the contents of row 21 (A21:D21) - the last data row (index 20,
# a Fibonacci-prime gap whose "gap" text was too large / errored out).
# Cell formatting (B21's style) is left in place, matching how Excel
# behaves when you select the cells and hit Delete.
$ws.Range("A21:D21").ClearContents()

# Excel moves the active cell up to B20 after the row-21 data is cleared.
$ws.Range("B20").Select()
